$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.589.22'
$ws.Range("E2").Value = '  -2.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.664.77'
$ws.Range("E3").Value = '  -3.65%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.63'
$ws.Range("E5").Value = '  -1.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.511'
$ws.Range("E6").Value = '  -2.20%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.93'
$ws.Range("E8").Value = '  -1.40%  '

$ws.Range("E9").Value = '  -0.83%  '

$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("E11").Value = '  -2.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.899.72'
$ws.Range("E12").Value = '  -3.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.654.66'
$ws.Range("E13").Value = '  -4.25%  '

$ws.Range("E14").Value = '  -3.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("E15").Value = '  +0.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.37'
$ws.Range("E16").Value = '  -1.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.565.94'
$ws.Range("E17").Value = '  -2.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.13'
$ws.Range("E18").Value = '  -0.57%  '

$ws.Range("E19").Value = '  -3.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.66'
$ws.Range("E20").Value = '  -4.65%  '

$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.49'
$ws.Range("E22").Value = '  -3.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.37'
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("E24").Value = '  -3.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.39'
$ws.Range("E25").Value = '  -1.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.23'
$ws.Range("E26").Value = '  -3.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.42'
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.112'
$ws.Range("E29").Value = '  -2.32%  '

$ws.Range("E30").Value = '  +2.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0503'
$ws.Range("E31").Value = '  -2.10%  '

$ws.Range("E32").Value = '  -2.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.458.90'
$ws.Range("E33").Value = '  -3.29%  '

$ws.Range("E34").Value = '  -4.71%  '

$ws.Range("E35").Value = '  -5.11%  '

$ws.Range("E36").Value = '  -1.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.925'
$ws.Range("E37").Value = '  -4.24%  '

$ws.Range("E38").Value = '  -1.41%  '

$ws.Range("E39").Value = '  -5.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").Value = '  -3.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.58'
$ws.Range("E41").Value = '  -1.92%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.22'
$ws.Range("E43").Value = '  -3.24%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.41'
$ws.Range("E44").Value = '  -5.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.807.54'
$ws.Range("E45").Value = '  -3.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.790'
$ws.Range("E46").Value = '  -2.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.73'
$ws.Range("E47").Value = '  -1.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.88'
$ws.Range("E48").Value = '  -2.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0108'
$ws.Range("E49").Value = '  -4.04%  '

$ws.Range("E50").Value = '  -2.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.92'
$ws.Range("E51").Value = '  -3.81%  '
